# Daily refresh of "days since previous payout" (col G) and
# "days until next payout" (col H->I) columns, simulating the
# passage of one calendar day (today moved from 2023-09-18 to 2023-09-19).
#
# For every data row:
#   - if column F (previous payout date) has a value, column G
#     (days since previous payout) is incremented by 1.
#   - if column H (next payout date) has a value, column I
#     (days to next payout) is decremented by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $fValue = $ws.Cells.Item($row, 6).Value()
    if ($fValue -ne $null -and $fValue -ne "") {
        $gCell = $ws.Cells.Item($row, 7)
        $gValue = $gCell.Value()
        if ($gValue -ne $null -and $gValue -ne "") {
            $gCell.Value = $gValue + 1
        }
    }

    $hValue = $ws.Cells.Item($row, 8).Value()
    if ($hValue -ne $null -and $hValue -ne "") {
        $iCell = $ws.Cells.Item($row, 9)
        $iValue = $iCell.Value()
        if ($iValue -ne $null -and $iValue -ne "") {
            $iCell.Value = $iValue - 1
        }
    }
}
